$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Copy the formatting of the last existing BOM row (row 20) down onto the
# new row so the new row picks up the same cell styles/borders used by
# the rest of the table.
$ws.Range("A20:K20").Copy()
$ws.Range("A21:K21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Add a new BOM line for the speaker on row 21 (row 20 was the last used row)
$ws.Range("A21").Value = ""
$ws.Range("B21").Value = "Speaker"
$ws.Range("C21").Value = "speaker"
$ws.Range("D21").Value = "S1"
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = "C49246995"
$ws.Range("H21").Value = ""
$ws.Range("I21").Value = ""
$ws.Range("J21").Value = ""
$ws.Range("K21").Value = ""

# Update selection to mirror the saved state in the workbook
$ws.Range("G21").Select()
